$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.501.80"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "1.620.98"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.29"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.523"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.03"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("E9").Value = "  +1.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0610"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0880"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("D12").Value = "1.851.58"
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("D13").Value = "1.615.76"
$ws.Range("E13").Value = "  -2.08%  "
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.549"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.27"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("D17").Value = "27.485.47"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.62"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.53"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.41%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.40"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.59%  "
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.12"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +8.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.64"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.87"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.51"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.17"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0483"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("E32").Value = "  -1.34%  "
$ws.Range("D33").Value = "1.464.75"
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("E34").Value = "  -2.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.54"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.948"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.77%  "
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.869"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.551"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.34%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  -2.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "67.16"
$ws.Range("D43").ClearFormats()
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.20"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.16%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.31"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.90%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.761.11"
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("E47").Value = "  +2.28%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.22"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0105"
$ws.Range("E49").Value = "  -1.62%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0994"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0503"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.14%  "
